# Update "Categoría Edad" sheet (2nd sheet) text labels in column C and D
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categoría Edad")

# Column C: "Educación Básica-Media (6 a 20 años)" -> "Educación Básica-Media (6 a 15 años)"
$ws.Range("C3").Value = "Educación Básica-Media (6 a 15 años)"
$ws.Range("C4").Value = "Educación Básica-Media (6 a 15 años)"

# Column D: "Inactivo" -> "Inactivo (0 a 15 años y más de 65 años)"
$ws.Range("D2").Value = "Inactivo (0 a 15 años y más de 65 años)"
$ws.Range("D3").Value = "Inactivo (0 a 15 años y más de 65 años)"
$ws.Range("D4").Value = "Inactivo (0 a 15 años y más de 65 años)"
$ws.Range("D15").Value = "Inactivo (0 a 15 años y más de 65 años)"
$ws.Range("D16").Value = "Inactivo (0 a 15 años y más de 65 años)"
$ws.Range("D17").Value = "Inactivo (0 a 15 años y más de 65 años)"
$ws.Range("D18").Value = "Inactivo (0 a 15 años y más de 65 años)"

# Column D: "Económicamente Activo" -> "Económicamente Activo (16 a 65 años)"
$ws.Range("D5").Value = "Económicamente Activo (16 a 65 años)"
$ws.Range("D6").Value = "Económicamente Activo (16 a 65 años)"
$ws.Range("D7").Value = "Económicamente Activo (16 a 65 años)"
$ws.Range("D8").Value = "Económicamente Activo (16 a 65 años)"
$ws.Range("D9").Value = "Económicamente Activo (16 a 65 años)"
$ws.Range("D10").Value = "Económicamente Activo (16 a 65 años)"
$ws.Range("D11").Value = "Económicamente Activo (16 a 65 años)"
$ws.Range("D12").Value = "Económicamente Activo (16 a 65 años)"
$ws.Range("D13").Value = "Económicamente Activo (16 a 65 años)"
$ws.Range("D14").Value = "Económicamente Activo (16 a 65 años)"

# Update the selection on this sheet to G11 (as recorded in the saved view)
$ws.Range("G11").Select()
